$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196, shifting existing rows 196:229 down to 197:230
$ws.Rows.Item(196).Insert(-4121)

# Populate the newly inserted row 196 with its data
$ws.Range("A196").Value = 5
$ws.Range("B196").Value = "Macroferia Regional de Talca"
$ws.Range("C196").Value = "Maule"
$ws.Range("D196").Value = 44522
$ws.Range("E196").Value = 7
$ws.Range("F196").Value = 100112023
$ws.Range("G196").Value = "Brócoli"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 5000
$ws.Range("K196").Value = 500
$ws.Range("L196").Value = 500
$ws.Range("M196").Value = 500
$ws.Range("N196").Value = "$/unidad"
$ws.Range("O196").Value = "Región del Maule"
$ws.Range("P196").Value = 500
$ws.Range("Q196").Value = 1
$ws.Range("R196").Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D
$ws.Range("D196").NumberFormat = $ws.Range("D197").NumberFormat
